# banco_dados_funcionarios.xlsx update
# - replace the single "data admissao" date column with separate
#   "ano" / "mes" / "dia" columns
# - add a "Bonus Excel" formula column that weights the bonus by time
#   of employment
# - adjust a few column widths / the active selection to match the
#   edited workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- headers -----------------------------------------------------
# F1 used to hold "data admissao"; it becomes "ano". G1/H1 are new.
$ws.Range("F1").Value = "ano"
$ws.Range("G1").Value = "mes"
$ws.Range("H1").Value = "dia"
$ws.Range("F1:H1").HorizontalAlignment = 1

$ws.Range("J1").Value = "Bonus Excel"

# ---- row 2 (Victor Wilson) ----------------------------------------
$ws.Range("F2").Value = 2012
$ws.Range("F2").NumberFormat = "0.00"

$ws.Range("G2").Value = 1
$ws.Range("G2").NumberFormat = "#,##0.00"

$ws.Range("H2").Value = 5
$ws.Range("H2").NumberFormat = "#,##0.00"

$ws.Range("I2").NumberFormat = "#,##0.00"

$ws.Range("J2").Formula = "=(((E2*3)+(E2*1))/(E2*5))*12"

$ws.Range("K2").NumberFormat = "yyyy-mm-dd"

# ---- row 3 (Flossie Wilson) ----------------------------------------
$ws.Range("F3").Value = 2012
$ws.Range("F3").NumberFormat = "0.00"

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 6

$ws.Range("K3").NumberFormat = "yyyy-mm-dd"

# ---- row 4 (Sherman Hodges) ----------------------------------------
$ws.Range("F4").Value = 2012
$ws.Range("F4").NumberFormat = "0.00"

$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 7

$ws.Range("K4").NumberFormat = "yyyy-mm-dd"

# ---- cosmetics -------------------------------------------------------
$ws.Columns(3).ColumnWidth = 27
$ws.Columns(7).ColumnWidth = 9.5
$ws.Columns(9).ColumnWidth = 10.17
$ws.Columns(10).ColumnWidth = 11.17
$ws.Columns(11).ColumnWidth = 9.5

$ws.Range("B10").Select()
